$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (AD, AE, AF) to row 1 with values,
# copying the existing header formatting (bold/bordered/centered) from
# the adjacent "Unnamed: 28" header cell (AC1) so the new headers match
# the look of the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row (2-57)
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 78
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 1
}

Write-Output "Season record columns added"
